$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 94, shifting existing rows 94-103 down to 95-104
$ws.Rows("94:94").Insert()

# Populate the newly inserted row 94 with the new weekly record
$ws.Range("A94").Value = 9
$ws.Range("B94").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C94").Value = "Metropolitana"
$ws.Range("D94").Value = 45077
$ws.Range("E94").Value = 13
$ws.Range("F94").Value = 100112029
$ws.Range("G94").Value = "Orégano"
$ws.Range("H94").Value = "Sin especificar"
$ws.Range("I94").Value = "Primera"
$ws.Range("J94").Value = 16
$ws.Range("K94").Value = 16000
$ws.Range("L94").Value = 18000
$ws.Range("M94").Value = 17000
$ws.Range("N94").Value = "$/docena de atados"
$ws.Range("O94").Value = "Región Metropolitana"
$ws.Range("P94").Value = 5667
$ws.Range("Q94").Value = 3
$ws.Range("R94").Value = "Hortaliza"
